# Timesheet "Submitted List" update: replace sample rows 2-4 with the full
# employee roster (31 rows) and add the send-email-reminder target rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds numbers-as-text ("1", "2", ...). Format as Text first so
# Excel stores the digits verbatim instead of converting them to numeric
# cells, then clear the format back off so no residual style sticks.
$idCol = $ws.Range('A2:A32')
$idCol.NumberFormat = '@'

$ws.Range('A2').Value = '1'
$ws.Range('B2').Value = 'Harin Vimal Bharathi.J'
$ws.Range('C2').Value = 'Yuvaraj Shanmugam'
$ws.Range('D2').Value = 'harinj@techcedence.com'

$ws.Range('A3').Value = '2'
$ws.Range('B3').Value = 'Sri Manikandan K'
$ws.Range('C3').Value = 'Poorna chandran R'
$ws.Range('D3').Value = 'srimanikandank@techcedence.com'

$ws.Range('A4').Value = '3'
$ws.Range('B4').Value = 'Krishna Kumar P'
$ws.Range('C4').Value = 'Karthick N Super Admin'
$ws.Range('D4').Value = 'krishnak@techcedence.com'

$ws.Range('A5').Value = '4'
$ws.Range('B5').Value = 'Harish S'
$ws.Range('C5').Value = 'harishs@techcedence.com'

$ws.Range('A6').Value = '5'
$ws.Range('B6').Value = 'Anjali Krishna'
$ws.Range('C6').Value = 'anjalik@techcedence.com'

$ws.Range('A7').Value = '6'
$ws.Range('B7').Value = 'Santhosh Kumar M'
$ws.Range('C7').Value = 'Prathamesh Rajput'
$ws.Range('D7').Value = 'santhoshk@techcedence.com'

$ws.Range('A8').Value = '7'
$ws.Range('B8').Value = 'Anil Kumar P'
$ws.Range('C8').Value = 'Srinivasan N'
$ws.Range('D8').Value = 'anilk@techcedence.com'

$ws.Range('A9').Value = '8'
$ws.Range('B9').Value = 'Manikandan R'
$ws.Range('C9').Value = 'Krishna Kumar P'
$ws.Range('D9').Value = 'manikandanr@techcedence.com'

$ws.Range('A10').Value = '9'
$ws.Range('B10').Value = 'Alan Singh P'
$ws.Range('C10').Value = 'Yugendran G'
$ws.Range('D10').Value = 'alans@techcedence.com'

$ws.Range('A11').Value = '10'
$ws.Range('B11').Value = 'Sathish Kumar'
$ws.Range('C11').Value = 'sathishk@techcedence.com'

$ws.Range('A12').Value = '11'
$ws.Range('B12').Value = 'Vivekkumar Perinbaraj'
$ws.Range('C12').Value = 'Saibali Barooah'
$ws.Range('D12').Value = 'Vivekkumarp@techcedence.com'

$ws.Range('A13').Value = '12'
$ws.Range('B13').Value = 'Dhinakaran VP'
$ws.Range('C13').Value = 'Venkatesh A'
$ws.Range('D13').Value = 'dhinakaranv@techcedence.com'

$ws.Range('A14').Value = '13'
$ws.Range('B14').Value = 'Venkateshwaran C'
$ws.Range('C14').Value = 'Venkateshwaranc@techcedence.com'

$ws.Range('A15').Value = '14'
$ws.Range('B15').Value = 'Niket Shah'
$ws.Range('C15').Value = 'nikets@techcedence.com'

$ws.Range('A16').Value = '15'
$ws.Range('B16').Value = 'Kathiravan P'
$ws.Range('C16').Value = 'Golda G'
$ws.Range('D16').Value = 'kathiravanp@techcedence.com'

$ws.Range('A17').Value = '16'
$ws.Range('B17').Value = 'Karthik Vinod'
$ws.Range('C17').Value = 'Prathamesh Rajput'
$ws.Range('D17').Value = 'karthikv@techcedence.com'

$ws.Range('A18').Value = '17'
$ws.Range('B18').Value = 'Saibali Barooah'
$ws.Range('C18').Value = 'Krishna Kumar P'
$ws.Range('D18').Value = 'saibalib@techcedence.com'

$ws.Range('A19').Value = '18'
$ws.Range('B19').Value = 'Yugendran G'
$ws.Range('C19').Value = 'Venkatesh A'
$ws.Range('D19').Value = 'yugendrang@techcedence.com'

$ws.Range('A20').Value = '19'
$ws.Range('B20').Value = 'Jayanth Kandregula'
$ws.Range('C20').Value = 'Venkatesh A'
$ws.Range('D20').Value = 'jayanthk@techcedence.com'

$ws.Range('A21').Value = '20'
$ws.Range('B21').Value = 'Gobi J'
$ws.Range('C21').Value = 'Krishna Kumar P'
$ws.Range('D21').Value = 'gobij@techcedence.com'

$ws.Range('A22').Value = '21'
$ws.Range('B22').Value = 'Manoowranjith A J'
$ws.Range('C22').Value = 'Saravana Kumar'
$ws.Range('D22').Value = 'Manoowranjitha@techcedence.com'

$ws.Range('A23').Value = '22'
$ws.Range('B23').Value = 'Sabarish K'
$ws.Range('C23').Value = 'Gowtham R'
$ws.Range('D23').Value = 'sabarishk@techcedence.com'

$ws.Range('A24').Value = '23'
$ws.Range('B24').Value = 'Saravana Kumar'
$ws.Range('C24').Value = 'Venkatesh A'
$ws.Range('D24').Value = 'saravanak@techcedence.com'

$ws.Range('A25').Value = '24'
$ws.Range('B25').Value = 'Prathamesh Rajput'
$ws.Range('C25').Value = 'Krishna Kumar P'
$ws.Range('D25').Value = 'prathameshr@techcedence.com'

$ws.Range('A26').Value = '25'
$ws.Range('B26').Value = 'Gowrishankar.G'
$ws.Range('C26').Value = 'Prathamesh Rajput'
$ws.Range('D26').Value = 'GowrishankarG@techcedence.com'

$ws.Range('A27').Value = '26'
$ws.Range('B27').Value = 'Venkat BCG'
$ws.Range('C27').Value = 'venkat@barcodegulf.net'

$ws.Range('A28').Value = '27'
$ws.Range('B28').Value = 'Mehdi S'
$ws.Range('C28').Value = 'mehdis@techcedence.com'

$ws.Range('A29').Value = '28'
$ws.Range('B29').Value = 'dina001'
$ws.Range('C29').Value = 'dhinakaranv+1@techcedence.com'

$ws.Range('A30').Value = '29'
$ws.Range('B30').Value = 'admin@techcedence.com'
$ws.Range('C30').Value = 'admin@techcedence.com'

$ws.Range('A31').Value = '30'
$ws.Range('B31').Value = 'Medi'
$ws.Range('C31').Value = 'mehdi.s@geevida.com'

$ws.Range('A32').Value = '31'
$ws.Range('B32').Value = 'US Cricket Store'
$ws.Range('C32').Value = 'uscricstore@gmail.com'

# Drop the temporary Text number-format now that the values are locked in
# as strings, so the cells end up with the default (no explicit) style.
$idCol.ClearFormats()

